$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new columns (D:L) before the old D column; old D:R shifts right to M:AA
$ws.Range("D1:L11").EntireColumn.Insert()

# Force text format on the numeric-looking new columns (K: fractionDigits, L: centAmount)
# so values are stored as text, matching the source data (numberStoredAsText).
$ws.Range("K1:L11").NumberFormat = "@"

# Header row (row 1) for the newly inserted columns
$ws.Range("D1").Value = 'variants.sku'
$ws.Range("E1").Value = 'description.en-US'
$ws.Range("F1").Value = 'priceMode'
$ws.Range("G1").Value = 'variants.prices.key'
$ws.Range("H1").Value = 'variants.prices.value.currencyCode'
$ws.Range("I1").Value = 'variants.prices.value.type'
$ws.Range("J1").Value = 'variants.prices.country'
$ws.Range("K1").Value = 'variants.prices.value.fractionDigits'
$ws.Range("L1").Value = 'variants.prices.value.centAmount'

# Data rows 2-11 for the newly inserted columns
# Row 2
$ws.Range("D2").Value = 'SNW-3125-01'
$ws.Range("E2").Value = 'The Lib Tech T.Rice Pro HP C2 Snowboard was born and raised in Jackson Hole. If that doesn''t mean anything to you, let us break it down: it''s annoyingly good at everything, loves a good cliff or couloir, and it invades Moab in the off season. If it can handle Corbet''s, it can definitely handle your neck of the woods, so level up your game with Travis Rice''s baby.'
$ws.Range("F2").Value = 'Embedded'
$ws.Range("G2").Value = 'SNW-Price-0-01'
$ws.Range("H2").Value = 'USD'
$ws.Range("I2").Value = 'centPrecision'
$ws.Range("J2").Value = 'US'
$ws.Range("K2").Value = '2'
$ws.Range("L2").Value = '48999'
# Row 3
$ws.Range("D3").Value = 'SNW-2365-01'
$ws.Range("E3").Value = 'Double down on your porcine persuasion with the Ride Twinpig Snowboard. This freestyle machine is everything you love about the short, fat, directional Warpig, but in a twin shape that means you get all the pork going forward, plus all the pork going backwards. It slays in either direction and does extra damage in the park with a wide shape and double blunted nose for easy spinning and stomping.'
$ws.Range("F3").Value = 'Embedded'
$ws.Range("G3").Value = 'SNW-Price-1-01'
$ws.Range("H3").Value = 'USD'
$ws.Range("I3").Value = 'centPrecision'
$ws.Range("J3").Value = 'US'
$ws.Range("K3").Value = '2'
$ws.Range("L3").Value = '37096'
# Row 4
$ws.Range("D4").Value = 'SNW-2537-01'
$ws.Range("E4").Value = 'Blast off into the shred-o-sphere with the CAPiTA Aeronaut Snowboard, an eye-popping cambered masterpiece thats''s made for the all-mountain grind. Created with legend Arthur Longo, it features a progressive sidecut and medium flex for versatility that''ll take you from side hits to the slackcountry. A lightweight and durable Hoover Core with Carbon Flax amplifiers add lightweight stability and tons of pop. The CAPiTA Aeronaut Snowboard is a jack of all trades, straight from the mothership.'
$ws.Range("F4").Value = 'Embedded'
$ws.Range("G4").Value = 'SNW-Price-2-01'
$ws.Range("H4").Value = 'USD'
$ws.Range("I4").Value = 'centPrecision'
$ws.Range("J4").Value = 'US'
$ws.Range("K4").Value = '2'
$ws.Range("L4").Value = '64995'
# Row 5
$ws.Range("D5").Value = 'SNW-2155-01'
$ws.Range("E5").Value = 'The Season Nexus Snowboard stands out in the arms race of bloated quivers and technical jargon - a simple skeleton key capable of unlocking the innate potential of any day in the mountains. Its straightforward geometry cuts through the noise, creating something that adds up to much more than the sum of its parts. Ample float in powder, smooth, engaged turns on hardpack, and dependable performance in everyday resort chop. It’s the ultimate quiver of one - more than enough to remind you why you fell in love with snowboarding in the first place.'
$ws.Range("F5").Value = 'Embedded'
$ws.Range("G5").Value = 'SNW-Price-3-01'
$ws.Range("H5").Value = 'USD'
$ws.Range("I5").Value = 'centPrecision'
$ws.Range("J5").Value = 'US'
$ws.Range("K5").Value = '2'
$ws.Range("L5").Value = '32940'
# Row 6
$ws.Range("D6").Value = 'SNW-2876-01'
$ws.Range("E6").Value = 'The Lib Tech T.Rice Orca Snowboard is still out there chomping fish. The volume-shifted pow ripper that''ll have you laughing as you chop the tops off moguls, the Orca has become the pillar of do-it-all boards. The sidecut''s tight enough to drill a spiral carve down to the bedrock, the short nose floats like a dang cork in the Dead Sea, and it stays stable at ludicrous speeds. From Jackson Hole to the BC BC, from Mt Baker to whatever that ice moon of Jupiter is called: the Orca''s ready to eat.'
$ws.Range("F6").Value = 'Embedded'
$ws.Range("G6").Value = 'SNW-Price-4-01'
$ws.Range("H6").Value = 'USD'
$ws.Range("I6").Value = 'centPrecision'
$ws.Range("J6").Value = 'US'
$ws.Range("K6").Value = '2'
$ws.Range("L6").Value = '48999'
# Row 7
$ws.Range("D7").Value = 'SNW-3398-01'
$ws.Range("E7").Value = 'The Rossignol Myth Snowboard and its reputation for making the transition from beginner to advanced a snap isn''t just some old wive''s tale, it''s the hard cold truth. The combination of Auto Turn Rocker and a relatively soft flex mean getting up and running has never been simpler. Toss out those falsehoods about snowboarding being cold, hard and difficult to master: get on the Myth and get after it.'
$ws.Range("F7").Value = 'Embedded'
$ws.Range("G7").Value = 'SNW-Price-5-01'
$ws.Range("H7").Value = 'USD'
$ws.Range("I7").Value = 'centPrecision'
$ws.Range("J7").Value = 'US'
$ws.Range("K7").Value = '2'
$ws.Range("L7").Value = '25997'
# Row 8
$ws.Range("D8").Value = 'SNW-2245-01'
$ws.Range("E8").Value = 'At this point, the Ride Warpig Snowboard is the stuff of legend. A card carrying member of the volume shifted revolution, the war-mongering hog sports a sturdy directional shape, a super wide platform, and a nimble dual-radius sidecut, delivering all mountain capability with a signature playful personality like few other boards can. The Warpig is one of the best all-arounders on the market, and it''s high time you found out what all the fuss is about.'
$ws.Range("F8").Value = 'Embedded'
$ws.Range("G8").Value = 'SNW-Price-6-01'
$ws.Range("H8").Value = 'USD'
$ws.Range("I8").Value = 'centPrecision'
$ws.Range("J8").Value = 'US'
$ws.Range("K8").Value = '2'
$ws.Range("L8").Value = '38496'
# Row 9
$ws.Range("D9").Value = 'SNW-2905-01'
$ws.Range("E9").Value = 'Sip it and rip it: the Lib Tech Cold Brew C2 Snowboard is pure frozen nitro go-go juice. An all-mountain ripper with a directional shape and powerful C2 rocker profile, this board is equally adept at quick shots down chutes, and thirsty-two ouncer long haul trucker runs. It does prefer fresh snow, but don''t we all? Order up a triple double with cream on top: this one''s the all day pick-me-up! Get low, stay powerful, remain caffeinated.'
$ws.Range("F9").Value = 'Embedded'
$ws.Range("G9").Value = 'SNW-Price-7-01'
$ws.Range("H9").Value = 'USD'
$ws.Range("I9").Value = 'centPrecision'
$ws.Range("J9").Value = 'US'
$ws.Range("K9").Value = '2'
$ws.Range("L9").Value = '38499'
# Row 10
$ws.Range("D10").Value = 'SNW-3162-01'
$ws.Range("E10").Value = 'Indulge in the smooth and effortless ride of the K2 Dreamsicle Snowboard, an all-mountain board designed to take you from silky side hits to soft stashes with ease. Combining a setback stance, twin shape, and directional rocker profile, the Dreamsicle offers precise control and catch-free turns, allowing you to carve easily on hardpack and float effortlessly on powder days. And with a supple flex and women''s-specific Rhythm Core, it''s the perfect board for relaxed cruising in all conditions.'
$ws.Range("F10").Value = 'Embedded'
$ws.Range("G10").Value = 'SNW-Price-8-01'
$ws.Range("H10").Value = 'USD'
$ws.Range("I10").Value = 'centPrecision'
$ws.Range("J10").Value = 'US'
$ws.Range("K10").Value = '2'
$ws.Range("L10").Value = '31496'
# Row 11
$ws.Range("D11").Value = 'SNW-2858-01'
$ws.Range("E11").Value = 'The CAPiTA Mega Mercury Snowboard takes the all-round versatility of the already popular Mercury and slathers it with extra sauce. Laid up with a super poppy core, this capable deck is amped up and slimmed down with Capita''s Megacarbon™ Array and Megalite Skin™ for lightweight performance that allows you take flight in any situation. New for this season are freshly minted wide models and Capita''s fastest Megadrive™ Base, found only on the Mega line-up.'
$ws.Range("F11").Value = 'Embedded'
$ws.Range("G11").Value = 'SNW-Price-9-01'
$ws.Range("H11").Value = 'USD'
$ws.Range("I11").Value = 'centPrecision'
$ws.Range("J11").Value = 'US'
$ws.Range("K11").Value = '2'
$ws.Range("L11").Value = '55996'

# Reset style on the text-forced columns back to Normal so no stray style index lingers,
# while the stored values remain text (shared-string) as already written above.
$ws.Range("K1:L11").Style = "Normal"

# Dimension / ignoredErrors sqref reflect the new A1:AA11 used range automatically via UsedRange;
# re-assert the numberStoredAsText ignored-error coverage over the full new grid.
